$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("roll")

# ---------------------------------------------------------------------------
# FRONT table (rows 3-13): ARB displacement doubled.
#   H = (D-F)*2   (was D-F)
#   I = (H-H_next)/(B-B_next)  -- now a shared formula block I3:I13
# ---------------------------------------------------------------------------
$ws.Range("H3").Formula = "=(D3-F3)*2"
$ws.Range("H4:H13").Formula = "=(D4-F4)*2"
$ws.Range("I3:I13").Formula = "=(H3-H4)/(B3-B4)"

# ---------------------------------------------------------------------------
# REAR table (rows 17-27): formulas unchanged, I becomes a shared formula
# block (I17:I27) instead of individual formulas - values stay the same.
# ---------------------------------------------------------------------------
$ws.Range("I17:I27").Formula = "=(H17-H18)/(B17-B18)"

# ---------------------------------------------------------------------------
# New reference data block K19:P29 (values only, no formulas)
# ---------------------------------------------------------------------------
$kpData = @(
  @(0, 2.5,  0, -7.2885,              242.07, 161.06),
  @(0, 2,    0, -7.2392000000000003,  242.81, 160.82),
  @(0, 1.5,  0, -7.1897000000000002,  243.47, 160.58000000000001),
  @(0, 1,    0, -7.14,                244.06, 160.37),
  @(0, 0.5,  0, -7.0900999999999996,  244.57, 160.16999999999999),
  @(0, 0,    0, -7.04,                245,    160),
  @(0, -0.5, 0, -6.9898999999999996,  245.35, 159.85),
  @(0, -1,   0, -6.9397000000000002,  245.61, 159.74),
  @(0, -1.5, 0, -6.8895999999999997,  245.79, 159.66),
  @(0, -2,   0, -6.8395000000000001,  245.86, 159.63),
  @(0, -2.5, 0, -6.7896999999999998,  245.82, 159.65)
)

$r = 19
foreach ($row in $kpData) {
    $c = 11  # column K
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# sortState now covers the new block
# ---------------------------------------------------------------------------
$ws.Sort.SetRange($ws.Range("K19:P29"))
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("L19"), 0, 2)
$ws.Sort.Header = 0
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# Selection moves to J9
# ---------------------------------------------------------------------------
$ws.Range("J9").Select()
